# Generate Report for Handoff
# Marks additional e2e files as high-priority ("ht") handoff candidates and
# refreshes the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
# timestamps for the files that were (re)handed off.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Priority column ("E") updates: flag rows as "ht" -----------------------

# zh-cn: all of rows 7-13 become "ht"
foreach ($row in 7..13) {
    $wsZhCn.Range("E$row").Value = "ht"
}

# de-de: rows 7, 9, 10, 11, 12 become "ht" (8 and 13 stay as-is)
foreach ($row in @(7, 9, 10, 11, 12)) {
    $wsDeDe.Range("E$row").Value = "ht"
}

# --- Latest Handoff Datetime refresh ----------------------------------------

# zh-cn "Latest Handoff Datetime" (column H) for rows 7-13: 10:25:29 -> 10:25:55
foreach ($row in 7..13) {
    $wsZhCn.Range("H$row").Value = "2016-08-18 10:25:55"
}

# de-de "Latest Handoff Datetime" (column H) for rows 7, 9, 10, 11, 12: 10:25:35 -> 10:26:04
foreach ($row in @(7, 9, 10, 11, 12)) {
    $wsDeDe.Range("H$row").Value = "2016-08-18 10:26:04"
}

# Overview "Latest HO Xliff Generate Date" (column G) mirrors the same refresh
$wsOverview.Range("G7").Value  = "2016-08-18 10:26:04"
$wsOverview.Range("G8").Value  = "2016-08-18 10:25:55"
$wsOverview.Range("G9").Value  = "2016-08-18 10:26:04"
$wsOverview.Range("G10").Value = "2016-08-18 10:26:04"
$wsOverview.Range("G11").Value = "2016-08-18 10:26:04"
$wsOverview.Range("G12").Value = "2016-08-18 10:26:04"
$wsOverview.Range("G13").Value = "2016-08-18 10:25:55"
